# Update the "Portugal Primeira Liga" results table.
#
# Two pairs of rows were out of order (match id in column B should be
# ascending) - row 14/15 and row 296/297 have their entire match record
# (everything except the running index in column A) swapped back into the
# correct order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 14 and 15: swap columns B:AD (keep A, the row index, fixed) ---
$range14 = $ws.Range("B14:AD14")
$range15 = $ws.Range("B15:AD15")

$vals14 = $range14.Value2
$vals15 = $range15.Value2

$range14.Value2 = $vals15
$range15.Value2 = $vals14

# --- Rows 296 and 297: swap columns B:AD (keep A, the row index, fixed) ---
$range296 = $ws.Range("B296:AD296")
$range297 = $ws.Range("B297:AD297")

$vals296 = $range296.Value2
$vals297 = $range297.Value2

$range296.Value2 = $vals297
$range297.Value2 = $vals296
